$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data value in C2 (30 -> 15)
$ws.Range("C2").Value = 15

# Move the active selection to F4 (was K11)
[void]$ws.Range("F4").Select()
